# Update cryptos price/volume data (and a few reordered rows) per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '29.079.52'
$ws.Cells.Item(2, 5).Value = '  -0.69%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.902.00'
$ws.Cells.Item(3, 5).Value = '  -0.65%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.001'
$ws.Cells.Item(4, 5).Value = '  -0.15%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '327.36'
$ws.Cells.Item(5, 5).Value = '  -0.06%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.9987'
$ws.Cells.Item(6, 5).Value = '  -0.42%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4607'
$ws.Cells.Item(7, 5).Value = '  -0.49%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3878'
$ws.Cells.Item(8, 5).Value = '  -1.75%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07852'
$ws.Cells.Item(9, 5).Value = '  -1.15%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.9886'
$ws.Cells.Item(10, 5).Value = '  -1.32%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '21.94'
$ws.Cells.Item(11, 5).Value = '  -1.87%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.864.78'
$ws.Cells.Item(12, 5).Value = '  -2.73%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.745'
$ws.Cells.Item(13, 5).Value = '  -0.49%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '7.025'
$ws.Cells.Item(14, 5).Value = '  -1.24%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.07015'
$ws.Cells.Item(15, 5).Value = '  +1.07%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '87.94'
$ws.Cells.Item(16, 5).Value = '  -0.78%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '1.002'
$ws.Cells.Item(17, 5).Value = '  -0.04%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000009925'
$ws.Cells.Item(18, 5).Value = '  -1.39%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '17.05'
$ws.Cells.Item(19, 5).Value = '  -0.41%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.9992'
$ws.Cells.Item(20, 5).Value = '  -0.29%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '29.085.62'
$ws.Cells.Item(21, 5).Value = '  -0.74%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.317'
$ws.Cells.Item(22, 5).Value = '  -0.81%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '11.09'
$ws.Cells.Item(23, 5).Value = '  -0.23%  '

# Row 24
$ws.Cells.Item(24, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(24, 4).Value = '2.100.35'
$ws.Cells.Item(24, 5).Value = '  -2.54%  '

# Row 25
$ws.Cells.Item(25, 2).Value = 'Toncoin'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.089'
$ws.Cells.Item(25, 5).Value = '  +1.64%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '156.10'
$ws.Cells.Item(26, 5).Value = '  -0.60%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '19.42'
$ws.Cells.Item(27, 5).Value = '  -0.28%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '5.891'
$ws.Cells.Item(28, 5).Value = '  -5.04%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '118.62'
$ws.Cells.Item(29, 5).Value = '  -0.26%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.871'
$ws.Cells.Item(30, 5).Value = '  -6.27%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.09339'
$ws.Cells.Item(31, 5).Value = '  -0.59%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.8935'
$ws.Cells.Item(32, 5).Value = '  -3.52%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.221'
$ws.Cells.Item(33, 5).Value = '  -2.62%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.317'
$ws.Cells.Item(34, 5).Value = '  -3.03%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '3.134'
$ws.Cells.Item(35, 5).Value = '  -4.29%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.05775'
$ws.Cells.Item(36, 5).Value = '  -1.15%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.167'
$ws.Cells.Item(37, 5).Value = '  -2.80%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.02081'
$ws.Cells.Item(38, 5).Value = '  -1.41%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.9982'
$ws.Cells.Item(39, 5).Value = '  -0.39%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.5682'
$ws.Cells.Item(40, 5).Value = '  -1.27%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '7.636'
$ws.Cells.Item(41, 5).Value = '  -4.30%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.1804'
$ws.Cells.Item(42, 5).Value = '  +0.04%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '9.699'
$ws.Cells.Item(43, 5).Value = '  -2.62%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '11.88'
$ws.Cells.Item(44, 5).Value = '  -1.18%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.5341'
$ws.Cells.Item(45, 5).Value = '  -1.58%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'PEPE'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.000002754'
$ws.Cells.Item(46, 5).Value = '  +68.18%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'RenderToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.164'
$ws.Cells.Item(47, 5).Value = '  -6.52%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Cronos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.06975'
$ws.Cells.Item(48, 5).Value = '  -1.41%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'NEARProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.835'
$ws.Cells.Item(49, 5).Value = '  -2.47%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '112.90'
$ws.Cells.Item(50, 5).Value = '  -0.59%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'MXToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.525'
$ws.Cells.Item(51, 5).Value = '  -1.40%  '
